$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 5483.6665
$ws.Range("J10").Value = 8000
$ws.Range("L10").Value = 8000
$ws.Range("N10").Value = -8586
$ws.Range("H19").Value = 771.96
$ws.Range("I19").Value = 447.1111
$ws.Range("K19").Value = 447.1111
$ws.Range("M19").Value = -272.1111
$ws.Range("H33").Value = 340.29413
$ws.Range("I33").Value = 281
$ws.Range("J33").Value = 425
$ws.Range("K33").Value = 281
$ws.Range("L33").Value = 425
$ws.Range("M33").Value = -52
$ws.Range("N33").Value = -883
$ws.Range("H40").Value = 2052.4075
$ws.Range("I40").Value = 2094
$ws.Range("J40").Value = 1991.909
$ws.Range("K40").Value = 2094
$ws.Range("L40").Value = 1991.909
$ws.Range("M40").Value = -1919
$ws.Range("N40").Value = -2341.909
$ws.Range("H98").Value = 1154.4231
$ws.Range("I98").Value = 1125.625
$ws.Range("K98").Value = 1125.625
$ws.Range("M98").Value = 372.375
$ws.Range("H106").Value = 2540
$ws.Range("I106").Value = 1480
$ws.Range("J106").Value = 2937.5
$ws.Range("K106").Value = 1480
$ws.Range("L106").Value = 2937.5
$ws.Range("M106").Value = -849
$ws.Range("N106").Value = -4199.5
$ws.Range("H122").Value = 1154.4231
$ws.Range("I122").Value = 1125.625
$ws.Range("K122").Value = 3376.875
$ws.Range("M122").Value = -926.875
$ws.Range("H132").Value = 1906.54
$ws.Range("I132").Value = 1488.8975
$ws.Range("K132").Value = 4466.6925
$ws.Range("M132").Value = -1936.6925

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 6913
$ws.Range("I3").Value = 7188.3335
$ws.Range("J3").Value = 6500
$ws.Range("K3").Value = 7188.3335
$ws.Range("L3").Value = 6500
$ws.Range("M3").Value = -7073.3335
$ws.Range("N3").Value = -6730
$ws.Range("H5").Value = 94.2
$ws.Range("I5").Value = 57
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 57
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = -374
$ws.Range("H9").Value = 32454.5
$ws.Range("J9").Value = 32454.5
$ws.Range("L9").Value = 32454.5
$ws.Range("N9").Value = -32794.5
$ws.Range("H20").Value = 32454.5
$ws.Range("J20").Value = 32454.5
$ws.Range("L20").Value = 32454.5
$ws.Range("N20").Value = -32994.5
$ws.Range("H61").Value = 5088.1143
$ws.Range("I61").Value = 4132.773
$ws.Range("J61").Value = 6704.846
$ws.Range("K61").Value = 4132.773
$ws.Range("L61").Value = 6704.846
$ws.Range("M61").Value = -3920.773
$ws.Range("N61").Value = -7128.846
$ws.Range("H74").Value = 4985.737
$ws.Range("I74").Value = 2813.3076
$ws.Range("J74").Value = 9692.666999999999
$ws.Range("K74").Value = 2813.3076
$ws.Range("L74").Value = 9692.666999999999
$ws.Range("M74").Value = -1939.3076
$ws.Range("N74").Value = -11440.667
$ws.Range("H77").Value = 4985.737
$ws.Range("I77").Value = 2813.3076
$ws.Range("J77").Value = 9692.666999999999
$ws.Range("K77").Value = 14066.538
$ws.Range("L77").Value = 48463.335
$ws.Range("M77").Value = -9698.538
$ws.Range("N77").Value = -57199.335
$ws.Range("H132").Value = 3707.349
$ws.Range("I132").Value = 1187
$ws.Range("J132").Value = 6148.9375
$ws.Range("K132").Value = 3561
$ws.Range("L132").Value = 18446.8125
$ws.Range("M132").Value = -1031
$ws.Range("N132").Value = -23506.8125
$ws.Range("H136").Value = 5088.1143
$ws.Range("I136").Value = 4132.773
$ws.Range("J136").Value = 6704.846
$ws.Range("K136").Value = 12398.319
$ws.Range("L136").Value = 20114.538
$ws.Range("M136").Value = -9848.319
$ws.Range("N136").Value = -25214.538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 94.2
$ws.Range("I4").Value = 57
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 57
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = 58
$ws.Range("N4").Value = -380
$ws.Range("H20").Value = 1087.9412
$ws.Range("I20").Value = 1089.6364
$ws.Range("J20").Value = 1084.8334
$ws.Range("K20").Value = 1089.6364
$ws.Range("L20").Value = 1084.8334
$ws.Range("M20").Value = -842.6364000000001
$ws.Range("N20").Value = -1578.8334
$ws.Range("H94").Value = 1234.5714
$ws.Range("I94").Value = 1054.4706
$ws.Range("K94").Value = 1054.4706
$ws.Range("M94").Value = -603.4706000000001
$ws.Range("H134").Value = 2463.5066
$ws.Range("I134").Value = 2550.7812
$ws.Range("J134").Value = 1955.7273
$ws.Range("K134").Value = 7652.3436
$ws.Range("L134").Value = 5867.1819
$ws.Range("M134").Value = -5117.3436
$ws.Range("N134").Value = -10937.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 158.27272
$ws.Range("I22").Value = 127.28571
$ws.Range("K22").Value = 127.28571
$ws.Range("M22").Value = 222.71429
$ws.Range("H31").Value = 1675.3286
$ws.Range("I31").Value = 1278.9661
$ws.Range("J31").Value = 3801.2727
$ws.Range("K31").Value = 1278.9661
$ws.Range("L31").Value = 3801.2727
$ws.Range("M31").Value = -983.9661000000001
$ws.Range("N31").Value = -4391.2727
$ws.Range("H34").Value = 1675.3286
$ws.Range("I34").Value = 1278.9661
$ws.Range("J34").Value = 3801.2727
$ws.Range("K34").Value = 1278.9661
$ws.Range("L34").Value = 3801.2727
$ws.Range("M34").Value = -1076.9661
$ws.Range("N34").Value = -4205.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1518.6052
$ws.Range("J107").Value = 1886.1724
$ws.Range("L107").Value = 5658.5172
$ws.Range("N107").Value = -9498.5172
$ws.Range("H122").Value = 634.37933
$ws.Range("I122").Value = 547
$ws.Range("J122").Value = 696.05884
$ws.Range("K122").Value = 4923
$ws.Range("L122").Value = 6264.52956
$ws.Range("M122").Value = -2473
$ws.Range("N122").Value = -11164.52956
$ws.Range("H125").Value = 2502.4
$ws.Range("J125").Value = 2917
$ws.Range("L125").Value = 8751
$ws.Range("N125").Value = -18591

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5212.857
$ws.Range("I70").Value = 4669.5186
$ws.Range("J70").Value = 5718.724
$ws.Range("K70").Value = 4669.5186
$ws.Range("L70").Value = 5718.724
$ws.Range("M70").Value = -4399.5186
$ws.Range("N70").Value = -6258.724
$ws.Range("H73").Value = 5212.857
$ws.Range("I73").Value = 4669.5186
$ws.Range("J73").Value = 5718.724
$ws.Range("K73").Value = 4669.5186
$ws.Range("L73").Value = 5718.724
$ws.Range("M73").Value = -3733.5186
$ws.Range("N73").Value = -7590.724
$ws.Range("H80").Value = 7995.6
$ws.Range("I80").Value = 14237.5
$ws.Range("J80").Value = 3834.3333
$ws.Range("K80").Value = 14237.5
$ws.Range("L80").Value = 3834.3333
$ws.Range("M80").Value = -13239.5
$ws.Range("N80").Value = -5830.3333
$ws.Range("H82").Value = 40328
$ws.Range("J82").Value = 40328
$ws.Range("L82").Value = 40328
$ws.Range("N82").Value = -41094
$ws.Range("H83").Value = 7995.6
$ws.Range("I83").Value = 14237.5
$ws.Range("J83").Value = 3834.3333
$ws.Range("K83").Value = 71187.5
$ws.Range("L83").Value = 19171.6665
$ws.Range("M83").Value = -66195.5
$ws.Range("N83").Value = -29155.6665
$ws.Range("H85").Value = 40328
$ws.Range("J85").Value = 40328
$ws.Range("L85").Value = 40328
$ws.Range("N85").Value = -42980
$ws.Range("H86").Value = 20143
$ws.Range("J86").Value = 20143
$ws.Range("L86").Value = 20143
$ws.Range("N86").Value = -22515
$ws.Range("H88").Value = 31753.334
$ws.Range("J88").Value = 31753.334
$ws.Range("L88").Value = 31753.334
$ws.Range("N88").Value = -32655.334
$ws.Range("H89").Value = 20143
$ws.Range("J89").Value = 20143
$ws.Range("L89").Value = 60429
$ws.Range("N89").Value = -72285
$ws.Range("H91").Value = 31753.334
$ws.Range("J91").Value = 31753.334
$ws.Range("L91").Value = 31753.334
$ws.Range("N91").Value = -34873.334
$ws.Range("H122").Value = 8955.714
$ws.Range("I122").Value = 11740
$ws.Range("J122").Value = 1995
$ws.Range("K122").Value = 35220
$ws.Range("L122").Value = 5985
$ws.Range("M122").Value = -32770
$ws.Range("N122").Value = -10885
$ws.Range("H132").Value = 2851.2812
$ws.Range("I132").Value = 2750.889
$ws.Range("J132").Value = 2980.3572
$ws.Range("K132").Value = 8252.667000000001
$ws.Range("L132").Value = 8941.071599999999
$ws.Range("M132").Value = -5722.667000000001
$ws.Range("N132").Value = -14001.0716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 757.05884
$ws.Range("J22").Value = 826.6667
$ws.Range("L22").Value = 826.6667
$ws.Range("N22").Value = -1416.6667
$ws.Range("H27").Value = 757.05884
$ws.Range("J27").Value = 826.6667
$ws.Range("L27").Value = 826.6667
$ws.Range("N27").Value = -1040.6667
$ws.Range("H100").Value = 3432.9614
$ws.Range("I100").Value = 2091.0625
$ws.Range("J100").Value = 5580
$ws.Range("K100").Value = 2091.0625
$ws.Range("L100").Value = 5580
$ws.Range("M100").Value = -1550.0625
$ws.Range("N100").Value = -6662

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 11333.333
$ws.Range("J24").Value = 11333.333
$ws.Range("L24").Value = 11333.333
$ws.Range("N24").Value = -11793.333
$ws.Range("H122").Value = 3551.375
$ws.Range("J122").Value = 4557.615
$ws.Range("L122").Value = 13672.845
$ws.Range("N122").Value = -18572.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M5").Value = -821.6667
$ws.Range("N5").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
